$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.341.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.279.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '499.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.39'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.68%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0953'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.336'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.72'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.686.51'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.80'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.325.53'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000129'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.271.41'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.24'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '303.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.150'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.31'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.17'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.61'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.81%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0687'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.93'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.78'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +10.25%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.373'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.40'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.37'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.93%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '125.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.10%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.83'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0494'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0894'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.548'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '240.53'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.372'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0205'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.76'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.32'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.65'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.19%  '
